$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 662.34
$ws.Range("I15").Value = 662.34
$ws.Range("K15").Value = 1987.02
$ws.Range("M15").Value = -1818.02

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1587.1
$ws.Range("I45").Value = 1395.3334
$ws.Range("J45").Value = 1874.75
$ws.Range("K45").Value = 1395.3334
$ws.Range("L45").Value = 1874.75
$ws.Range("M45").Value = -1018.3334
$ws.Range("N45").Value = -2628.75
# Row 63
$ws.Range("H63").Value = 11349.333
$ws.Range("I63").Value = 15250
$ws.Range("J63").Value = 3548
$ws.Range("K63").Value = 15250
$ws.Range("L63").Value = 3548
$ws.Range("M63").Value = -14564
$ws.Range("N63").Value = -4920
# Row 66
$ws.Range("H66").Value = 11349.333
$ws.Range("I66").Value = 15250
$ws.Range("J66").Value = 3548
$ws.Range("K66").Value = 76250
$ws.Range("L66").Value = 17740
$ws.Range("M66").Value = -72818
$ws.Range("N66").Value = -24604
# Row 74
$ws.Range("H74").Value = 2272.6086
$ws.Range("I74").Value = 1301.9333
$ws.Range("J74").Value = 4092.625
$ws.Range("K74").Value = 1301.9333
$ws.Range("L74").Value = 4092.625
$ws.Range("M74").Value = -427.9332999999999
$ws.Range("N74").Value = -5840.625
# Row 77
$ws.Range("H77").Value = 2272.6086
$ws.Range("I77").Value = 1301.9333
$ws.Range("J77").Value = 4092.625
$ws.Range("K77").Value = 6509.666499999999
$ws.Range("L77").Value = 20463.125
$ws.Range("M77").Value = -2141.666499999999
$ws.Range("N77").Value = -29199.125
# Row 122
$ws.Range("H122").Value = 2009.3235
$ws.Range("I122").Value = 1565.5834
$ws.Range("J122").Value = 3074.3
$ws.Range("K122").Value = 4696.7502
$ws.Range("L122").Value = 9222.900000000001
$ws.Range("M122").Value = -2246.7502
$ws.Range("N122").Value = -14122.9

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1262.9231
$ws.Range("I7").Value = 584.75
$ws.Range("K7").Value = 584.75
$ws.Range("M7").Value = -471.75
# Row 59
$ws.Range("H59").Value = 27890
$ws.Range("J59").Value = 27890
$ws.Range("L59").Value = 27890
$ws.Range("N59").Value = -29584
# Row 134
$ws.Range("H134").Value = 24798.639
$ws.Range("J134").Value = 6128.5713
$ws.Range("L134").Value = 18385.7139
$ws.Range("N134").Value = -23455.7139

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("K12").Value = 300
$ws.Range("M12").Value = -130
# Row 122
$ws.Range("H122").Value = 1107.7391
$ws.Range("J122").Value = 1577.6666
$ws.Range("L122").Value = 4732.9998
$ws.Range("N122").Value = -9632.9998
# Row 132
$ws.Range("H132").Value = 2421.5881
$ws.Range("I132").Value = 1395
$ws.Range("J132").Value = 3334.111
$ws.Range("K132").Value = 4185
$ws.Range("L132").Value = 10002.333
$ws.Range("M132").Value = -1655
$ws.Range("N132").Value = -15062.333
# Row 134
$ws.Range("H134").Value = 2212.889
$ws.Range("I134").Value = 1345.4
$ws.Range("J134").Value = 2723.1765
$ws.Range("K134").Value = 4036.2
$ws.Range("L134").Value = 8169.529500000001
$ws.Range("M134").Value = -1501.2
$ws.Range("N134").Value = -13239.5295

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 13.05
$ws.Range("I12").Value = 20
$ws.Range("J12").Value = 8.416667
$ws.Range("K12").Value = 60
$ws.Range("L12").Value = 25.250001
$ws.Range("M12").Value = 113
$ws.Range("N12").Value = -371.250001
# Row 49
$ws.Range("H49").Value = 3003
$ws.Range("I49").Value = 3000
$ws.Range("J49").Value = 3004
$ws.Range("K49").Value = 9000
$ws.Range("L49").Value = 9012
$ws.Range("M49").Value = -8844
$ws.Range("N49").Value = -9324
# Row 131
$ws.Range("H131").Value = 796.4737
$ws.Range("I131").Value = 255.15384
$ws.Range("J131").Value = 908.1746000000001
$ws.Range("K131").Value = 765.4615200000001
$ws.Range("L131").Value = 2724.5238
$ws.Range("M131").Value = 4274.53848
$ws.Range("N131").Value = -12804.5238
# Row 132
$ws.Range("H132").Value = 6194.4546
$ws.Range("I132").Value = 829.875
$ws.Range("J132").Value = 20500
$ws.Range("K132").Value = 7468.875
$ws.Range("L132").Value = 184500
$ws.Range("M132").Value = -4938.875
$ws.Range("N132").Value = -189560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 9
$ws.Range("H9").Value = 755
$ws.Range("I9").Value = 755
$ws.Range("K9").Value = 755
$ws.Range("M9").Value = -585
# Row 31
$ws.Range("H31").Value = 377
$ws.Range("I31").Value = 377
$ws.Range("K31").Value = 377
$ws.Range("M31").Value = -85
# Row 37
$ws.Range("H37").Value = 377
$ws.Range("I37").Value = 377
$ws.Range("K37").Value = 377
$ws.Range("M37").Value = -100
# Row 70
$ws.Range("H70").Value = 11530.529
$ws.Range("I70").Value = 5067.933
$ws.Range("K70").Value = 5067.933
$ws.Range("M70").Value = -4797.933
# Row 73
$ws.Range("H73").Value = 11530.529
$ws.Range("I73").Value = 5067.933
$ws.Range("K73").Value = 5067.933
$ws.Range("M73").Value = -4131.933
# Row 102
$ws.Range("H102").Value = 2755.3333
$ws.Range("I102").Value = 2835.3635
$ws.Range("J102").Value = 1875
$ws.Range("K102").Value = 2835.3635
$ws.Range("L102").Value = 1875
$ws.Range("M102").Value = -1213.3635
$ws.Range("N102").Value = -5119
# Row 122
$ws.Range("H122").Value = 1280.5714
$ws.Range("I122").Value = 1155.4
$ws.Range("J122").Value = 1593.5
$ws.Range("K122").Value = 3466.2
$ws.Range("L122").Value = 4780.5
$ws.Range("M122").Value = -1016.2
$ws.Range("N122").Value = -9680.5
# Row 132
$ws.Range("H132").Value = 4941.0347
$ws.Range("I132").Value = 6172.2856
$ws.Range("J132").Value = 3791.8667
$ws.Range("K132").Value = 18516.8568
$ws.Range("L132").Value = 11375.6001
$ws.Range("M132").Value = -15986.8568
$ws.Range("N132").Value = -16435.6001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 318.16666
$ws.Range("I9").Value = 318.16666
$ws.Range("K9").Value = 318.16666
$ws.Range("M9").Value = -94.16665999999998
# Row 35
$ws.Range("H35").Value = 2029
$ws.Range("I35").Value = 2029
$ws.Range("K35").Value = 2029
$ws.Range("M35").Value = -1693
# Row 109
$ws.Range("H109").Value = 30663.334
$ws.Range("J109").Value = 30663.334
$ws.Range("L109").Value = 30663.334
$ws.Range("N109").Value = -33437.334
# Row 136
$ws.Range("H136").Value = 5369.35
$ws.Range("I136").Value = 3194.9473
$ws.Range("J136").Value = 7336.6665
$ws.Range("K136").Value = 9584.841899999999
$ws.Range("L136").Value = 22009.9995
$ws.Range("M136").Value = -7034.841899999999
$ws.Range("N136").Value = -27109.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 57360.832
$ws.Range("I122").Value = 126364.625
$ws.Range("K122").Value = 379093.875
$ws.Range("M122").Value = -376643.875
# Row 132
$ws.Range("H132").Value = 25812.232
$ws.Range("I132").Value = 84891.586
$ws.Range("K132").Value = 254674.758
$ws.Range("M132").Value = -252144.758
# Row 133
$ws.Range("H133").Value = 29665
$ws.Range("J133").Value = 29665
$ws.Range("L133").Value = 29665
$ws.Range("N133").Value = -39785
